# Documentation fix for the KLayout PCell scripting guide title slide:
# "Klayout: Parameterized Cells Scriptnig" -> "Klayout: Parameterized Cells Scripting"
# (moves "Cells " into the run containing the fixed word, matching the
#  original author's retyping of the misspelled word).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)   # "Title 11" placeholder shape
$tr = $shp.TextFrame.TextRange

$fullText = $tr.Text
$oldFragment = "Cells Scriptnig"
$newFragment = "Cells Scripting"

$startPos = $fullText.IndexOf($oldFragment)
if ($startPos -ge 0) {
    # Characters() is 1-based.
    $target = $tr.Characters($startPos + 1, $oldFragment.Length)
    $target.Text = $newFragment
}
